# Auto-generated edit script: updates currentAveragePrice / Leve price & profit
# columns (H-N) for specific rows across multiple crafting-leve sheets, to
# reflect a refreshed pull from the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 894.5
$ws.Range("I6").Value = 437.66666
$ws.Range("J6").Value = 1808.1666
$ws.Range("K6").Value = 1312.99998
$ws.Range("L6").Value = 5424.4998
$ws.Range("M6").Value = -1200.99998
$ws.Range("N6").Value = -5648.4998
$ws.Range("H12").Value = 1049.875
$ws.Range("I12").Value = 600
$ws.Range("J12").Value = 1499.75
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 1499.75
$ws.Range("M12").Value = -430
$ws.Range("N12").Value = -1839.75
$ws.Range("H40").Value = 2860.6
$ws.Range("I40").Value = 2729.8
$ws.Range("K40").Value = 2729.8
$ws.Range("M40").Value = -2554.8
$ws.Range("H64").Value = 3210
$ws.Range("I64").Value = 3012.5
$ws.Range("K64").Value = 3012.5
$ws.Range("M64").Value = -2764.5
$ws.Range("H67").Value = 3210
$ws.Range("I67").Value = 3012.5
$ws.Range("K67").Value = 3012.5
$ws.Range("M67").Value = -2154.5
$ws.Range("H74").Value = 3563.8572
$ws.Range("I74").Value = 2990.6
$ws.Range("K74").Value = 2990.6
$ws.Range("M74").Value = -2054.6
$ws.Range("H77").Value = 3563.8572
$ws.Range("I77").Value = 2990.6
$ws.Range("K77").Value = 14953
$ws.Range("M77").Value = -10273
$ws.Range("H137").Value = 56881.555
$ws.Range("I137").Value = 983
$ws.Range("J137").Value = 92453.37
$ws.Range("K137").Value = 2949
$ws.Range("L137").Value = 277360.11
$ws.Range("M137").Value = -399
$ws.Range("N137").Value = -282460.11

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5495.3013
$ws.Range("I32").Value = 3678.2876
$ws.Range("J32").Value = 16676.924
$ws.Range("K32").Value = 3678.2876
$ws.Range("L32").Value = 16676.924
$ws.Range("M32").Value = -3391.2876
$ws.Range("N32").Value = -17250.924
$ws.Range("H45").Value = 1567.3636
$ws.Range("I45").Value = 1124.25
$ws.Range("K45").Value = 1124.25
$ws.Range("M45").Value = -747.25
$ws.Range("H74").Value = 715.2432
$ws.Range("I74").Value = 722.8889
$ws.Range("K74").Value = 722.8889
$ws.Range("M74").Value = 151.1111
$ws.Range("H77").Value = 715.2432
$ws.Range("I77").Value = 722.8889
$ws.Range("K77").Value = 3614.4445
$ws.Range("M77").Value = 753.5554999999999
$ws.Range("H88").Value = 2465.2
$ws.Range("I88").Value = 1979.5454
$ws.Range("K88").Value = 1979.5454
$ws.Range("M88").Value = -1573.5454
$ws.Range("H91").Value = 2465.2
$ws.Range("I91").Value = 1979.5454
$ws.Range("K91").Value = 1979.5454
$ws.Range("M91").Value = -575.5454
$ws.Range("H97").Value = 1526.5
$ws.Range("I97").Value = 1380.7
$ws.Range("J97").Value = 2255.5
$ws.Range("K97").Value = 1380.7
$ws.Range("L97").Value = 2255.5
$ws.Range("M97").Value = -884.7
$ws.Range("N97").Value = -3247.5
$ws.Range("H132").Value = 1871.0182
$ws.Range("J132").Value = 2214.5386
$ws.Range("L132").Value = 6643.6158
$ws.Range("N132").Value = -11703.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 144126.92
$ws.Range("I86").Value = 922.375
$ws.Range("K86").Value = 922.375
$ws.Range("M86").Value = 200.625
$ws.Range("H89").Value = 144126.92
$ws.Range("I89").Value = 922.375
$ws.Range("K89").Value = 4611.875
$ws.Range("M89").Value = 1004.125
$ws.Range("H94").Value = 352.25
$ws.Range("I94").Value = 352.25
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 352.25
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 98.75
$ws.Range("N94").Value = ""
$ws.Range("H107").Value = 3115.5
$ws.Range("I107").Value = 2838.6
$ws.Range("K107").Value = 2838.6
$ws.Range("M107").Value = -918.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3492.6155
$ws.Range("I31").Value = 2275.8
$ws.Range("K31").Value = 2275.8
$ws.Range("M31").Value = -1980.8
$ws.Range("H34").Value = 3492.6155
$ws.Range("I34").Value = 2275.8
$ws.Range("K34").Value = 2275.8
$ws.Range("M34").Value = -2073.8
$ws.Range("H58").Value = 4349978
$ws.Range("I58").Value = 5437070.5
$ws.Range("J58").Value = 1607
$ws.Range("K58").Value = 5437070.5
$ws.Range("L58").Value = 1607
$ws.Range("M58").Value = -5436867.5
$ws.Range("N58").Value = -2013
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736
$ws.Range("H136").Value = 4349978
$ws.Range("I136").Value = 5437070.5
$ws.Range("J136").Value = 1607
$ws.Range("K136").Value = 16311211.5
$ws.Range("L136").Value = 4821
$ws.Range("M136").Value = -16308661.5
$ws.Range("N136").Value = -9921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
$ws.Range("H17").Value = 13217.5
$ws.Range("I17").Value = 149.5
$ws.Range("J17").Value = 19751.5
$ws.Range("K17").Value = 448.5
$ws.Range("L17").Value = 59254.5
$ws.Range("M17").Value = -279.5
$ws.Range("N17").Value = -59592.5
$ws.Range("H104").Value = 3661.5386
$ws.Range("I104").Value = 200
$ws.Range("J104").Value = 5200
$ws.Range("K104").Value = 600
$ws.Range("L104").Value = 15600
$ws.Range("M104").Value = 2021
$ws.Range("N104").Value = -20842
$ws.Range("H107").Value = 509.15384
$ws.Range("J107").Value = 615.5
$ws.Range("L107").Value = 1846.5
$ws.Range("N107").Value = -5686.5
$ws.Range("H131").Value = 22474.312
$ws.Range("J131").Value = 26562.408
$ws.Range("L131").Value = 79687.224
$ws.Range("N131").Value = -89767.224
$ws.Range("H139").Value = 5202.25
$ws.Range("I139").Value = 5618.0415
$ws.Range("J139").Value = 2707.5
$ws.Range("K139").Value = 16854.1245
$ws.Range("L139").Value = 8122.5
$ws.Range("M139").Value = -11714.1245
$ws.Range("N139").Value = -18402.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H80").Value = 2780
$ws.Range("I80").Value = 2784.2856
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 2784.2856
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -1786.2856
$ws.Range("N80").Value = -4746
$ws.Range("H83").Value = 2780
$ws.Range("I83").Value = 2784.2856
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 13921.428
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -8929.428
$ws.Range("N83").Value = -23734
$ws.Range("H122").Value = 2484.9
$ws.Range("I122").Value = 2449.889
$ws.Range("K122").Value = 7349.667
$ws.Range("M122").Value = -4899.667
$ws.Range("H132").Value = 1103533.4
$ws.Range("I132").Value = 1543651.2
$ws.Range("K132").Value = 4630953.6
$ws.Range("M132").Value = -4628423.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11727.706
$ws.Range("I40").Value = 11491.4
$ws.Range("K40").Value = 11491.4
$ws.Range("M40").Value = -11355.4
$ws.Range("H82").Value = 2363.8462
$ws.Range("I82").Value = 1901.25
$ws.Range("K82").Value = 1901.25
$ws.Range("M82").Value = -1540.25
$ws.Range("H85").Value = 2363.8462
$ws.Range("I85").Value = 1901.25
$ws.Range("K85").Value = 1901.25
$ws.Range("M85").Value = -653.25
$ws.Range("H122").Value = 2913.25
$ws.Range("I122").Value = 2761.4
$ws.Range("K122").Value = 8284.200000000001
$ws.Range("M122").Value = -5834.200000000001
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""
$ws.Range("H132").Value = 5101.8887
$ws.Range("I132").Value = 2764.2856
$ws.Range("J132").Value = 6589.4546
$ws.Range("K132").Value = 8292.856800000001
$ws.Range("L132").Value = 19768.3638
$ws.Range("M132").Value = -5762.856800000001
$ws.Range("N132").Value = -24828.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 90007
$ws.Range("J9").Value = 90007
$ws.Range("L9").Value = 90007
$ws.Range("N9").Value = -90287
$ws.Range("H122").Value = 68289.086
$ws.Range("I122").Value = 90396.55499999999
$ws.Range("K122").Value = 271189.665
$ws.Range("M122").Value = -268739.665
